$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 3's formatting (cell style + row height) down onto row 4,
# same as the source row ("passive income" / "...") that row 4 follows.
$ws.Range("A3:B3").Copy()
$null = $ws.Range("A4:B4").PasteSpecial(-4122)
$ws.Rows.Item(4).RowHeight = $ws.Rows.Item(3).RowHeight

# New keyword pair.
$ws.Range("A4").Value = "affiliate marketing"
$ws.Range("B4").Value = "affiliate.marketing.guide"

$null = $ws.Range("B4").Select()
